# Applies the pull-request edit described in the commit:
#  - adds a new trailing column AK (a plain 0 in the header row, matching the
#    header style already used by B1:AJ1)
#  - fixes the swapped tkm-N3Usage / tkm-SZMUsage column headers (and swaps
#    the two data values in row 2 to stay under the correct header)
#  - updates a batch of recomputed figures in row 2
#  - adds a new row 3 (index "2" in column A) carrying a new label,
#    "BASF Schwarzheide GmbH", in the new column AK

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new column AK in the header row ---------------------------------
$ws.Range("AK1").Value = 0
$ws.Range("B1").Copy()
$ws.Range("AK1").PasteSpecial(-4122)   # xlPasteFormats - reuse the bold/border header style

# ---- fix the mislabeled tkm-N3Usage / tkm-SZMUsage headers ------------
$ws.Range("AI1").Value = "tkm-N3Usage"
$ws.Range("AJ1").Value = "tkm-SZMUsage"

# ---- row 2: recomputed figures -----------------------------------------
$ws.Range("B2").Value = 60754.97540209295
$ws.Range("E2").Value = 400
$ws.Range("F2").Value = 1214701.412959892
$ws.Range("J2").Value = 216
$ws.Range("K2").Value = 216
$ws.Range("N2").Value = 999999
$ws.Range("O2").Value = 999999
$ws.Range("P2").Value = 20350.18688630918
$ws.Range("Q2").Value = 20350.18688630918
$ws.Range("R2").Value = 16354.37667599064
$ws.Range("S2").Value = 16354.37667599064
$ws.Range("T2").Value = 833653.6233240096
$ws.Range("U2").Value = 833653.6233240094
$ws.Range("AD2").Value = 20715.54378958815
$ws.Range("AE2").Value = 1000956.785227869
$ws.Range("AF2").Value = 850008

# values under the (now corrected) AI1/AJ1 headers swap places
$ws.Range("AI2").Value = 130.3
$ws.Range("AJ2").Value = 414.5

# ---- row 3: new row ------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)    # xlPasteFormats - reuse A2's style

$ws.Range("AK3").Value = "BASF Schwarzheide GmbH"
